$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Test Suite")

# Update Runmode column (C) for rows 3-7 from "Y" to "N"
$ws.Range("C3").Value = "N"
$ws.Range("C4").Value = "N"
$ws.Range("C5").Value = "N"
$ws.Range("C6").Value = "N"
$ws.Range("C7").Value = "N"

# Update selection to A6
$ws.Activate()
$ws.Range("A6").Select()
